$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = "example_ques.xlsx"

$ws.Range("B1").Select()
